$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly record was inserted before the old row 82, pushing the
# existing rows 82-89 down to 83-90 (data unchanged), and the new row 82
# is populated with the latest price report.
$ws.Rows("82:82").Insert()

$ws.Range("A82").Value = 10
$ws.Range("B82").Value = "Vega Modelo de Temuco"
$ws.Range("C82").Value = "La Araucanía"
$ws.Range("D82").Value = 45194
$ws.Range("E82").Value = 9
$ws.Range("F82").Value = 300000000
$ws.Range("G82").Value = "Espárragos"
$ws.Range("H82").Value = "Sin especificar"
$ws.Range("I82").Value = "Primera"
$ws.Range("J82").Value = 550
$ws.Range("K82").Value = 1700
$ws.Range("L82").Value = 2000
$ws.Range("M82").Value = 1836
$ws.Range("N82").Value = "`$/kilo"
$ws.Range("O82").Value = "Región de La Araucanía"
$ws.Range("P82").Value = 1836
$ws.Range("Q82").Value = 1
$ws.Range("R82").Value = "Hortaliza"
